$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "k-means(test)" sheet: update heading text + selection
# ---------------------------------------------------------------------
$wsTest = $wb.Worksheets.Item("k-means(test)")
$wsTest.Range("A10").Value = "compress ratio (iteration = 1000, clusters = 100)"
$wsTest.Range("F3").Select()

# ---------------------------------------------------------------------
# 2. "k-means(obs_info)" sheet: update heading texts + selection
#    (these results are now for the 2366316-element dataset)
# ---------------------------------------------------------------------
$wsObs = $wb.Worksheets.Item("k-means(obs_info)")
$wsObs.Range("A1").Value = "execution time (2366316, iteration = 1000, clusters = 100, procs = 4)"
$wsObs.Range("H2").Value = "c-0.1"
$wsObs.Range("A10").Value = "compress ratio (iteration = 1000, clusters = 100)"
$wsObs.Range("A15").Value = "execution time (2366316, iteration = 1000, clusters = 1000, procs = 16)"
$wsObs.Range("G3").Select()

# ---------------------------------------------------------------------
# 3. Add new "k-means(num_plasma)" sheet right after "k-means(obs_info)"
#    (and right before "fftss") with the num_plasma (4386200-element)
#    results, mirroring the obs_info layout.
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $wsObs)
$wsNew.Name = "k-means(num_plasma)"

$wsNew.Range("A1").Value = "execution time (4386200, iteration = 1000, clusters = 100, procs = 4)"

$wsNew.Range("A2").Value = "absErrBound"
$wsNew.Range("B2").Value = "uc"
$wsNew.Range("C2").Value = "c-0.000001"
$wsNew.Range("D2").Value = "c-0.00001"
$wsNew.Range("E2").Value = "c-0.0001"
$wsNew.Range("F2").Value = "c-0.001"
$wsNew.Range("G2").Value = "c-0.01"
$wsNew.Range("H2").Value = "c-0.1"

$wsNew.Range("B3").Value = 2206.6678849999998
$wsNew.Range("C3").Value = 2159.9150559999998
$wsNew.Range("D3").Value = 2220.9988400000002
$wsNew.Range("E3").Value = 2105.8497699999998
$wsNew.Range("F3").Value = 2205.4593639999998
$wsNew.Range("G3").Value = 2172.3178710000002
$wsNew.Range("H3").Value = 2165.6020990000002

$wsNew.Range("A5").Value = "real gosa"

$wsNew.Range("A6").Value = "absErrBound"
$wsNew.Range("B6").Value = 0.000001
$wsNew.Range("C6").Value = 0.00001
$wsNew.Range("D6").Value = 0.0001
$wsNew.Range("E6").Value = 0.001
$wsNew.Range("F6").Value = 0.01
$wsNew.Range("G6").Value = 0.1

$wsNew.Range("A7").Value = "compressed (sz-mod)"
$wsNew.Range("B7").Value = 0
$wsNew.Range("C7").Value = 0
$wsNew.Range("D7").Value = [double]"3.0000000000000001E-6"
$wsNew.Range("E7").Value = [double]"8.7000000000000001E-5"
$wsNew.Range("F7").Value = [double]"1.8320000000000001E-3"
$wsNew.Range("G7").Value = [double]"1.8320000000000001E-3"

$wsNew.Range("A10").Value = "compress ratio (iteration = 1000, clusters = 100)"

$wsNew.Range("A11").Value = "absErrBound"
$wsNew.Range("B11").Value = 0.000001
$wsNew.Range("C11").Value = 0.00001
$wsNew.Range("D11").Value = 0.0001
$wsNew.Range("E11").Value = 0.001
$wsNew.Range("F11").Value = 0.01
$wsNew.Range("G11").Value = 0.1

$wsNew.Range("B12").Value = 1.0000180000000001
$wsNew.Range("C12").Value = 1.0004420000000001
$wsNew.Range("D12").Value = 1.115305
$wsNew.Range("E12").Value = 1.7523679999999999
$wsNew.Range("F12").Value = 6.6116710000000003
$wsNew.Range("G12").Value = 6.6116710000000003

$wsNew.Range("A15").Value = "execution time (4386200, iteration = 1000, clusters = 1000, procs = 16)"

$wsNew.Range("A16").Value = "absErrBound"
$wsNew.Range("B16").Value = "uc"
$wsNew.Range("C16").Value = "c-0.000001"
$wsNew.Range("D16").Value = "c-0.00001"
$wsNew.Range("E16").Value = "c-0.0001"
$wsNew.Range("F16").Value = "c-0.001"
$wsNew.Range("G16").Value = "c-0.01"

$wsNew.Range("A19").Value = "real gosa"

$wsNew.Range("A20").Value = "absErrBound"
$wsNew.Range("B20").Value = 0.000001
$wsNew.Range("C20").Value = 0.00001
$wsNew.Range("D20").Value = 0.0001
$wsNew.Range("E20").Value = 0.001
$wsNew.Range("F20").Value = 0.01
$wsNew.Range("G20").Value = 0.1

$wsNew.Range("A21").Value = "compressed (sz-mod)"

$wsNew.Range("A24").Value = "compress ratio (iteration = 1000, clusters = 1000)"

$wsNew.Range("A25").Value = "absErrBound"
$wsNew.Range("B25").Value = 0.000001
$wsNew.Range("C25").Value = 0.00001
$wsNew.Range("D25").Value = 0.0001
$wsNew.Range("E25").Value = 0.001
$wsNew.Range("F25").Value = 0.01
$wsNew.Range("G25").Value = 0.1

$wsNew.Range("G3").Select()
